$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C32 (was "NA" text, becomes an empty text value)
$ws.Range("C32").Value = "'"
$ws.Range("C32").Style = "Normal"

# Add new row 33
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "2025-03-26"
$ws.Range("A33").Style = "Normal"
$ws.Range("B33").Value = "bonnes pratiques"
$ws.Range("C33").Value = 13
$ws.Range("D33").Value = 1
